$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 41682812
$ws.Range("I86").Value = 52651204
$ws.Range("J86").Value = 2920.8
$ws.Range("K86").Value = 52651204
$ws.Range("L86").Value = 2920.8
$ws.Range("M86").Value = -52650081
$ws.Range("N86").Value = -5166.8

$ws.Range("H89").Value = 41682812
$ws.Range("I89").Value = 52651204
$ws.Range("J89").Value = 2920.8
$ws.Range("K89").Value = 263256020
$ws.Range("L89").Value = 14604
$ws.Range("M89").Value = -263250404
$ws.Range("N89").Value = -25836

$ws.Range("H116").Value = 3273.6667
$ws.Range("I116").Value = 3255
$ws.Range("J116").Value = 3311
$ws.Range("K116").Value = 3255
$ws.Range("L116").Value = 3311
$ws.Range("M116").Value = 187
$ws.Range("N116").Value = -10195

$ws.Range("H132").Value = 1576.2051
$ws.Range("I132").Value = 1442.8182
$ws.Range("J132").Value = 2309.8333
$ws.Range("K132").Value = 4328.4546
$ws.Range("L132").Value = 6929.499899999999
$ws.Range("M132").Value = -1798.4546
$ws.Range("N132").Value = -11989.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1690.75
$ws.Range("I45").Value = 1232.4615
$ws.Range("J45").Value = 2541.8572
$ws.Range("K45").Value = 1232.4615
$ws.Range("L45").Value = 2541.8572
$ws.Range("M45").Value = -855.4614999999999
$ws.Range("N45").Value = -3295.8572

$ws.Range("H61").Value = 8550897
$ws.Range("I61").Value = 15153372
$ws.Range("J61").Value = 6517.647
$ws.Range("K61").Value = 15153372
$ws.Range("L61").Value = 6517.647
$ws.Range("M61").Value = -15153160
$ws.Range("N61").Value = -6941.647

$ws.Range("H122").Value = 78634.16
$ws.Range("I122").Value = 92567.63
$ws.Range("K122").Value = 277702.89
$ws.Range("M122").Value = -275252.89

$ws.Range("H136").Value = 8550897
$ws.Range("I136").Value = 15153372
$ws.Range("J136").Value = 6517.647
$ws.Range("K136").Value = 45460116
$ws.Range("L136").Value = 19552.941
$ws.Range("M136").Value = -45457566
$ws.Range("N136").Value = -24652.941

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 84923.836
$ws.Range("I107").Value = 126260.75
$ws.Range("J107").Value = 2250
$ws.Range("K107").Value = 126260.75
$ws.Range("L107").Value = 2250
$ws.Range("M107").Value = -124340.75
$ws.Range("N107").Value = -6090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1731.5625
$ws.Range("I16").Value = 2205.5
$ws.Range("J16").Value = 1447.2
$ws.Range("K16").Value = 2205.5
$ws.Range("L16").Value = 1447.2
$ws.Range("M16").Value = -1918.5
$ws.Range("N16").Value = -2021.2

$ws.Range("H23").Value = 33151.54
$ws.Range("I23").Value = 8498.333000000001
$ws.Range("J23").Value = 54282.855
$ws.Range("K23").Value = 8498.333000000001
$ws.Range("L23").Value = 54282.855
$ws.Range("M23").Value = -8258.333000000001
$ws.Range("N23").Value = -54762.855

$ws.Range("H27").Value = 33151.54
$ws.Range("I27").Value = 8498.333000000001
$ws.Range("J27").Value = 54282.855
$ws.Range("K27").Value = 8498.333000000001
$ws.Range("L27").Value = 54282.855
$ws.Range("M27").Value = -8306.333000000001
$ws.Range("N27").Value = -54666.855

$ws.Range("H31").Value = 4692.854
$ws.Range("I31").Value = 1257.4286
$ws.Range("J31").Value = 7762.8086
$ws.Range("K31").Value = 1257.4286
$ws.Range("L31").Value = 7762.8086
$ws.Range("M31").Value = -962.4286
$ws.Range("N31").Value = -8352.8086

$ws.Range("H34").Value = 4692.854
$ws.Range("I34").Value = 1257.4286
$ws.Range("J34").Value = 7762.8086
$ws.Range("K34").Value = 1257.4286
$ws.Range("L34").Value = 7762.8086
$ws.Range("M34").Value = -1055.4286
$ws.Range("N34").Value = -8166.8086

$ws.Range("H105").Value = 918
$ws.Range("I105").Value = 883.75
$ws.Range("J105").Value = 1055
$ws.Range("K105").Value = 883.75
$ws.Range("L105").Value = 1055
$ws.Range("M105").Value = 863.25
$ws.Range("N105").Value = -4549

$ws.Range("H113").Value = 1731.5625
$ws.Range("I113").Value = 2205.5
$ws.Range("J113").Value = 1447.2
$ws.Range("K113").Value = 2205.5
$ws.Range("L113").Value = 1447.2
$ws.Range("M113").Value = -35.5
$ws.Range("N113").Value = -5787.2

$ws.Range("H122").Value = 1916.4286
$ws.Range("I122").Value = 1575.8334
$ws.Range("J122").Value = 2052.6667
$ws.Range("K122").Value = 4727.5002
$ws.Range("L122").Value = 6158.000100000001
$ws.Range("M122").Value = -2277.5002
$ws.Range("N122").Value = -11058.0001

$ws.Range("H134").Value = 5226.3667
$ws.Range("I134").Value = 5166.0835
$ws.Range("K134").Value = 15498.2505
$ws.Range("M134").Value = -12963.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 73.454544
$ws.Range("J12").Value = 111.72727
$ws.Range("L12").Value = 335.18181
$ws.Range("N12").Value = -681.18181

$ws.Range("H121").Value = 1560.6666
$ws.Range("J121").Value = 1679.8182
$ws.Range("L121").Value = 5039.4546
$ws.Range("N121").Value = -7659.4546

$ws.Range("H131").Value = 2757.651
$ws.Range("J131").Value = 3414.1226
$ws.Range("L131").Value = 10242.3678
$ws.Range("N131").Value = -20322.3678

$ws.Range("H136").Value = 3655.3635
$ws.Range("I136").Value = 1845.4445
$ws.Range("J136").Value = 11800
$ws.Range("K136").Value = 5536.333500000001
$ws.Range("L136").Value = 35400
$ws.Range("M136").Value = -436.3335000000006
$ws.Range("N136").Value = -45600

$ws.Range("H137").Value = 32543.783
$ws.Range("I137").Value = 7009.35
$ws.Range("J137").Value = 62584.293
$ws.Range("K137").Value = 21028.05
$ws.Range("L137").Value = 187752.879
$ws.Range("M137").Value = -15928.05
$ws.Range("N137").Value = -197952.879

$ws.Range("H139").Value = 438647.47
$ws.Range("I139").Value = 910356.5600000001
$ws.Range("J139").Value = 6247.5
$ws.Range("K139").Value = 2731069.68
$ws.Range("L139").Value = 18742.5
$ws.Range("M139").Value = -2725929.68
$ws.Range("N139").Value = -29022.5

$ws.Range("H140").Value = 1059.3928
$ws.Range("I140").Value = 930.7406999999999
$ws.Range("K140").Value = 2792.2221
$ws.Range("M140").Value = 2387.7779

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 2541.4
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 2926.75
$ws.Range("K29").Value = 1000
$ws.Range("L29").Value = 2926.75
$ws.Range("M29").Value = -710
$ws.Range("N29").Value = -3506.75

$ws.Range("H122").Value = 8061.5625
$ws.Range("I122").Value = 10834.363
$ws.Range("K122").Value = 32503.089
$ws.Range("M122").Value = -30053.089

$ws.Range("H123").Value = 8214.038
$ws.Range("J123").Value = 8422.6
$ws.Range("L123").Value = 8422.6
$ws.Range("N123").Value = -13322.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 60173.332
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 90010
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 90010
$ws.Range("M4").Value = -387
$ws.Range("N4").Value = -90236

$ws.Range("H23").Value = 7004.6665
$ws.Range("J23").Value = 7004.6665
$ws.Range("L23").Value = 7004.6665
$ws.Range("M23").Value = -7464.6665

$ws.Range("H25").Value = 10000000
$ws.Range("J25").Value = 10000000
$ws.Range("L25").Value = 10000000
$ws.Range("N25").Value = -10000460

$ws.Range("H28").Value = 60173.332
$ws.Range("I28").Value = 500
$ws.Range("J28").Value = 90010
$ws.Range("K28").Value = 500
$ws.Range("L28").Value = 90010
$ws.Range("M28").Value = -268
$ws.Range("N28").Value = -90474

$ws.Range("H37").Value = 60173.332
$ws.Range("I37").Value = 500
$ws.Range("J37").Value = 90010
$ws.Range("K37").Value = 500
$ws.Range("L37").Value = 90010
$ws.Range("M37").Value = -393
$ws.Range("N37").Value = -90224

$ws.Range("H61").Value = 3660.4
$ws.Range("I61").Value = 2675.375
$ws.Range("K61").Value = 2675.375
$ws.Range("M61").Value = -2473.375

$ws.Range("H113").Value = 3660.4
$ws.Range("I113").Value = 2675.375
$ws.Range("K113").Value = 2675.375
$ws.Range("M113").Value = -505.375

$ws.Range("H122").Value = 5166.1113
$ws.Range("I122").Value = 5066.6665
$ws.Range("J122").Value = 5186
$ws.Range("K122").Value = 15199.9995
$ws.Range("L122").Value = 15558
$ws.Range("M122").Value = -12749.9995
$ws.Range("N122").Value = -20458

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1261.5555
$ws.Range("I113").Value = 1511.1666
$ws.Range("J113").Value = 762.3333
$ws.Range("K113").Value = 4533.4998
$ws.Range("L113").Value = 2286.9999
$ws.Range("M113").Value = -2363.4998
$ws.Range("N113").Value = -6626.9999

$ws.Range("H122").Value = 2084.1
$ws.Range("I122").Value = 2119.5557
$ws.Range("K122").Value = 6358.6671
$ws.Range("M122").Value = -3908.6671
